$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
Write-Host $ws.Name
